$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row 17
$ws.Range("H17").Value = 2916837.2
$ws.Range("J17").Value = 2970849.8
$ws.Range("L17").Value = 8912549.399999999
$ws.Range("N17").Value = -8912885.399999999

# ALC!row 20
$ws.Range("H20").Value = 6146.6665
$ws.Range("I20").Value = 2000
$ws.Range("J20").Value = 8220
$ws.Range("K20").Value = 2000
$ws.Range("L20").Value = 8220
$ws.Range("M20").Value = -1770
$ws.Range("N20").Value = -8680

# ALC!row 35
$ws.Range("H35").Value = 6146.6665
$ws.Range("I35").Value = 2000
$ws.Range("J35").Value = 8220
$ws.Range("K35").Value = 2000
$ws.Range("L35").Value = 8220
$ws.Range("M35").Value = -1621
$ws.Range("N35").Value = -8978

# ALC!row 51
$ws.Range("H51").Value = 2480
$ws.Range("I51").Value = 1900
$ws.Range("J51").Value = 3350
$ws.Range("K51").Value = 1900
$ws.Range("L51").Value = 3350
$ws.Range("M51").Value = -1416
$ws.Range("N51").Value = -4318

# ALC!row 112
$ws.Range("H112").Value = 11364706
$ws.Range("I112").Value = 621.25
$ws.Range("J112").Value = 13890058
$ws.Range("K112").Value = 1863.75
$ws.Range("L112").Value = 41670174
$ws.Range("M112").Value = -755.75
$ws.Range("N112").Value = -41672390

# ALC!row 132
$ws.Range("H132").Value = 2558.2856
$ws.Range("I132").Value = 2558.2856
$ws.Range("K132").Value = 7674.8568
$ws.Range("M132").Value = -5144.8568

# ALC!row 138
$ws.Range("H138").Value = 3391.5432
$ws.Range("I138").Value = 1373.1
$ws.Range("J138").Value = 4578.863
$ws.Range("K138").Value = 4119.299999999999
$ws.Range("L138").Value = 13736.589
$ws.Range("M138").Value = 1020.700000000001
$ws.Range("N138").Value = -24016.589

$ws = $wb.Worksheets.Item("ARM")
# ARM!row 5
$ws.Range("I5").Value = 300
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 300
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = -188
$ws.Range("N5").Value = -424

# ARM!row 52
$ws.Range("H52").Value = 35945
$ws.Range("J52").Value = 35945
$ws.Range("L52").Value = 35945
$ws.Range("N52").Value = -36581

# ARM!row 132
$ws.Range("H132").Value = 2099.9343
$ws.Range("I132").Value = 1169.9546
$ws.Range("J132").Value = 4506.9414
$ws.Range("K132").Value = 3509.8638
$ws.Range("L132").Value = 13520.8242
$ws.Range("M132").Value = -979.8638000000001
$ws.Range("N132").Value = -18580.8242

$ws = $wb.Worksheets.Item("BSM")
# BSM!row 4
$ws.Range("I4").Value = 300
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -185
$ws.Range("N4").Value = -430

$ws = $wb.Worksheets.Item("CRP")
# CRP!row 7
$ws.Range("H7").Value = 394.83334
$ws.Range("I7").Value = 394.83334
$ws.Range("K7").Value = 394.83334
$ws.Range("M7").Value = -281.83334

# CRP!row 31
$ws.Range("H31").Value = 2351.811
$ws.Range("I31").Value = 1907.2727
$ws.Range("J31").Value = 2495.6323
$ws.Range("K31").Value = 1907.2727
$ws.Range("L31").Value = 2495.6323
$ws.Range("M31").Value = -1612.2727
$ws.Range("N31").Value = -3085.6323

# CRP!row 34
$ws.Range("H34").Value = 2351.811
$ws.Range("I34").Value = 1907.2727
$ws.Range("J34").Value = 2495.6323
$ws.Range("K34").Value = 1907.2727
$ws.Range("L34").Value = 2495.6323
$ws.Range("M34").Value = -1705.2727
$ws.Range("N34").Value = -2899.6323

# CRP!row 64
$ws.Range("H64").Value = 25271
$ws.Range("J64").Value = 25271
$ws.Range("L64").Value = 25271
$ws.Range("N64").Value = -25767

# CRP!row 67
$ws.Range("H67").Value = 25271
$ws.Range("J67").Value = 25271
$ws.Range("L67").Value = 25271
$ws.Range("N67").Value = -26987

# CRP!row 69
$ws.Range("H69").Value = 25750
$ws.Range("I69").Value = 9633.333000000001
$ws.Range("J69").Value = 41866.668
$ws.Range("K69").Value = 9633.333000000001
$ws.Range("L69").Value = 41866.668
$ws.Range("M69").Value = -8884.333000000001
$ws.Range("N69").Value = -43364.668

# CRP!row 72
$ws.Range("H72").Value = 25750
$ws.Range("I72").Value = 9633.333000000001
$ws.Range("J72").Value = 41866.668
$ws.Range("K72").Value = 28899.999
$ws.Range("L72").Value = 125600.004
$ws.Range("M72").Value = -25155.999
$ws.Range("N72").Value = -133088.004

# CRP!row 74
$ws.Range("H74").Value = 28899.5
$ws.Range("J74").Value = 28899.5
$ws.Range("L74").Value = 28899.5
$ws.Range("N74").Value = -30647.5

# CRP!row 77
$ws.Range("H77").Value = 28899.5
$ws.Range("J77").Value = 28899.5
$ws.Range("L77").Value = 86698.5
$ws.Range("N77").Value = -95434.5

# CRP!row 134
$ws.Range("H134").Value = 5093.5386
$ws.Range("I134").Value = 5093.5386
$ws.Range("K134").Value = 15280.6158
$ws.Range("M134").Value = -12745.6158

# CRP!row 137
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

# CRP!row 140
$ws.Range("H140").Value = 21477.322
$ws.Range("J140").Value = 21477.322
$ws.Range("L140").Value = 21477.322
$ws.Range("N140").Value = -31837.322

# CRP!row 141
$ws.Range("H141").Value = 69396.336
$ws.Range("J141").Value = 69396.336
$ws.Range("L141").Value = 69396.336
$ws.Range("N141").Value = -79756.336

$ws = $wb.Worksheets.Item("CUL")
# CUL!row 68
$ws.Range("H68").Value = 2723.766
$ws.Range("I68").Value = 3498.7896
$ws.Range("J68").Value = 1968.6154
$ws.Range("K68").Value = 10496.3688
$ws.Range("L68").Value = 5905.8462
$ws.Range("M68").Value = -9685.3688
$ws.Range("N68").Value = -7527.8462

# CUL!row 71
$ws.Range("H71").Value = 2723.766
$ws.Range("I71").Value = 3498.7896
$ws.Range("J71").Value = 1968.6154
$ws.Range("K71").Value = 31489.1064
$ws.Range("L71").Value = 17717.5386
$ws.Range("M71").Value = -27433.1064
$ws.Range("N71").Value = -25829.5386

# CUL!row 92
$ws.Range("H92").Value = 661.5
$ws.Range("I92").Value = 725
$ws.Range("J92").Value = 534.5
$ws.Range("K92").Value = 2175
$ws.Range("L92").Value = 1603.5
$ws.Range("M92").Value = -927
$ws.Range("N92").Value = -4099.5

# CUL!row 107
$ws.Range("H107").Value = 1121.4286
$ws.Range("I107").Value = 324.84616
$ws.Range("J107").Value = 1409.0834
$ws.Range("K107").Value = 974.5384799999999
$ws.Range("L107").Value = 4227.2502
$ws.Range("M107").Value = 945.4615200000001
$ws.Range("N107").Value = -8067.2502

$ws = $wb.Worksheets.Item("LTW")
# LTW!row 111
$ws.Range("H111").Value = 50387
$ws.Range("J111").Value = 50387
$ws.Range("L111").Value = 50387
$ws.Range("N111").Value = -58567

# LTW!row 136
$ws.Range("H136").Value = 6650.49
$ws.Range("I136").Value = 5345.5312
$ws.Range("J136").Value = 8848.315000000001
$ws.Range("K136").Value = 16036.5936
$ws.Range("L136").Value = 26544.945
$ws.Range("M136").Value = -13486.5936
$ws.Range("N136").Value = -31644.945

# LTW!row 139
$ws.Range("H139").Value = 63900
$ws.Range("J139").Value = 63900
$ws.Range("L139").Value = 63900
$ws.Range("N139").Value = -74180

$ws = $wb.Worksheets.Item("WVR")
# WVR!row 123
$ws.Range("H123").Value = 24282
$ws.Range("J123").Value = 24282
$ws.Range("L123").Value = 24282
$ws.Range("N123").Value = -34082

# WVR!row 136
$ws.Range("H136").Value = 1467.6
$ws.Range("I136").Value = 869.4211
$ws.Range("K136").Value = 2608.2633
$ws.Range("M136").Value = -58.26330000000007
